$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "ec2"
$ws.Name = "ec2"

# Update the "type" column values from "vm" to "ec2"
$ws.Range("A2").Value = "ec2"
$ws.Range("A3").Value = "ec2"

# Remove the now-unused "tenancy" and "operatingsystem" columns (G:H)
$ws.Range("G1:H3").EntireColumn.Delete()

# Move the active selection to F2, matching the saved view state
$ws.Range("F2").Select()
